# Auto-generated script applying numeric updates to Sheet1 (pl_mw.xlsx case)
# per commit "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:C, rows 2:25
$arr = New-Object 'object[,]' 24,2
$arr[0,0] = 1.741765415654982
$arr[0,1] = 0.315201221576018
$arr[1,0] = 1.669426246988337
$arr[1,1] = 0.2802518922745776
$arr[2,0] = 1.626379214194287
$arr[2,1] = 0.2589322668135594
$arr[3,0] = 1.609180134560233
$arr[3,1] = 0.250278239534822
$arr[4,0] = 1.606344912034956
$arr[4,1] = 0.2488432583304814
$arr[5,0] = 1.62614587470938
$arr[5,1] = 0.2588154198671759
$arr[6,0] = 1.716537718865766
$arr[6,1] = 0.3031210462149261
$arr[7,0] = 1.904741616135482
$arr[7,1] = 0.3911674829723211
$arr[8,0] = 2.049818628007188
$arr[8,1] = 0.456650332556876
$arr[9,0] = 2.117325742166315
$arr[9,1] = 0.4866320641676793
$arr[10,0] = 2.143108395939009
$arr[10,1] = 0.4980146131021002
$arr[11,0] = 2.137545860997989
$arr[11,1] = 0.4955618611067507
$arr[12,0] = 2.11944249279054
$arr[12,1] = 0.4875679220606344
$arr[13,0] = 2.10838225815246
$arr[13,1] = 0.4826752354047699
$arr[14,0] = 2.045437443698518
$arr[14,1] = 0.4546949576536576
$arr[15,0] = 2.007211145061035
$arr[15,1] = 0.4375803982171647
$arr[16,0] = 1.985366478879826
$arr[16,1] = 0.4277546655048923
$arr[17,0] = 1.977994597092561
$arr[17,1] = 0.4244309114408225
$arr[18,0] = 2.011265680413828
$arr[18,1] = 0.4394003885055326
$arr[19,0] = 2.124753926896631
$arr[19,1] = 0.4899151347823931
$arr[20,0] = 2.200203515533588
$arr[20,1] = 0.5230997963093955
$arr[21,0] = 2.159817019110392
$arr[21,1] = 0.5053724854143979
$arr[22,0] = 2.00943221353117
$arr[22,1] = 0.4385775284625879
$arr[23,0] = 1.852642386014224
$arr[23,1] = 0.367215246460205
$ws.Range("B2:C25").Value = $arr

# Columns E:G, rows 2:25
$arr = New-Object 'object[,]' 24,3
$arr[0,0] = 0.03682411140076436
$arr[0,1] = 0.4443680307746263
$arr[0,2] = 0.002567512075434409
$arr[1,0] = 0.03611878584069572
$arr[1,1] = 0.387822817061874
$arr[1,2] = 0.002573702394190318
$arr[2,0] = 0.03567722831596054
$arr[2,1] = 0.3531389305168915
$arr[2,2] = 0.002577698083355261
$arr[3,0] = 0.03549512543775535
$arr[3,1] = 0.3390132514313251
$arr[3,2] = 0.002579375527679803
$arr[4,0] = 0.03546475575603303
$arr[4,1] = 0.336668177824194
$arr[4,2] = 0.002579657041010369
$arr[5,0] = 0.03567478121818102
$arr[5,1] = 0.3529483938368969
$arr[5,2] = 0.002577720506516368
$arr[6,0] = 0.03658265725436838
$arr[6,1] = 0.4248636149813336
$arr[6,2] = 0.002569606175233583
$arr[7,0] = 0.03829714036171517
$arr[7,1] = 0.5661985755042025
$arr[7,2] = 0.002555231354481575
$arr[8,0] = 0.03951874133826472
$arr[8,1] = 0.6702781546542269
$arr[8,2] = 0.002545595649030523
$arr[9,0] = 0.04006668980245109
$arr[9,1] = 0.7176906081379002
$arr[9,2] = 0.002541410558638214
$arr[10,0] = 0.04027310404765139
$arr[10,1] = 0.7356546913071611
$arr[10,2] = 0.002539854087190813
$arr[11,0] = 0.04022869661162964
$arr[11,1] = 0.7317853510981394
$arr[11,2] = 0.002540188044057506
$arr[12,0] = 0.0400836931108639
$arr[12,1] = 0.7191683204515869
$arr[12,2] = 0.002541281940069525
$arr[13,0] = 0.03999473441196777
$arr[13,1] = 0.7114413442032514
$arr[13,2] = 0.002541955667391745
$arr[14,0] = 0.03948277832386182
$arr[14,1] = 0.6671810134426437
$arr[14,2] = 0.002545873128469984
$arr[15,0] = 0.03916674568984213
$arr[15,1] = 0.6400460337215605
$arr[15,2] = 0.0025483270133507
$arr[16,0] = 0.03898423920187533
$arr[16,1] = 0.6244449056556647
$arr[16,2] = 0.002549757092843899
$arr[17,0] = 0.03892231872682395
$arr[17,1] = 0.6191636801734006
$arr[17,2] = 0.002550244504981647
$arr[18,0] = 0.03920046350492967
$arr[18,1] = 0.6429339538360921
$arr[18,2] = 0.002548063862054053
$arr[19,0] = 0.04012631324004268
$arr[19,1] = 0.7228739723492197
$arr[19,2] = 0.002540959868646907
$arr[20,0] = 0.04072511955859248
$arr[20,1] = 0.7751780083420101
$arr[20,2] = 0.002536482054163125
$arr[21,0] = 0.040406089350709
$arr[21,1] = 0.7472568307916134
$arr[21,2] = 0.002538856903442086
$arr[22,0] = 0.03918522222221377
$arr[22,1] = 0.6416283278902171
$arr[22,2] = 0.002548182772575586
$arr[23,0] = 0.0378401605944596
$arr[23,1] = 0.5279251897347308
$arr[23,2] = 0.002558956754530868
$ws.Range("E2:G25").Value = $arr

# Columns I:I, rows 2:25
$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 2.522390347086386
$arr[1,0] = 2.454662371584902
$arr[2,0] = 2.413784777052655
$arr[3,0] = 2.397301171849634
$arr[4,0] = 2.394574506805071
$arr[5,0] = 2.413561772214692
$arr[6,0] = 2.498888677896531
$arr[7,0] = 2.672000130297093
$arr[8,0] = 2.802960667245316
$arr[9,0] = 2.863412347724278
$arr[10,0] = 2.886434144293332
$arr[11,0] = 2.881470141167654
$arr[12,0] = 2.86530373559566
$arr[13,0] = 2.85541839997714
$arr[14,0] = 2.799027992407673
$arr[15,0] = 2.764661733880615
$arr[16,0] = 2.744977488359865
$arr[17,0] = 2.738326759745377
$arr[18,0] = 2.768311532824868
$arr[19,0] = 2.870048639563151
$arr[20,0] = 2.937299641339564
$arr[21,0] = 2.901335674843921
$arr[22,0] = 2.766661230474853
$arr[23,0] = 2.624523911064102
$ws.Range("I2:I25").Value = $arr

# Columns L:M, rows 2:25
$arr = New-Object 'object[,]' 24,2
$arr[0,0] = 0.2346688179144678
$arr[0,1] = 0.3657213672982849
$arr[1,0] = 0.2313587299036897
$arr[1,1] = 0.3540580275311029
$arr[2,0] = 0.2294434660531479
$arr[2,1] = 0.3471591268825378
$arr[3,0] = 0.2286923540941501
$arr[3,1] = 0.3444135255197196
$arr[4,0] = 0.2285694042408366
$arr[4,1] = 0.3439615853487581
$arr[5,0] = 0.2294332174691505
$arr[5,1] = 0.347121832798841
$arr[6,0] = 0.2335031341971145
$arr[6,1] = 0.3616452300354567
$arr[7,0] = 0.2424182312798564
$arr[7,1] = 0.3922209249390178
$arr[8,0] = 0.2495449277952986
$arr[8,1] = 0.4159838391813366
$arr[9,0] = 0.2529139922245918
$arr[9,1] = 0.4270812858245634
$arr[10,0] = 0.2542081679857091
$arr[10,1] = 0.4313253156086319
$arr[11,0] = 0.2539286246826435
$arr[11,1] = 0.4304094295536558
$arr[12,0] = 0.2530200957659048
$arr[12,1] = 0.4274296077873743
$arr[13,0] = 0.252465992999646
$arr[13,1] = 0.4256098167891977
$arr[14,0] = 0.2493273181489286
$arr[14,1] = 0.4152644083440791
$arr[15,0] = 0.2474344706712657
$arr[15,1] = 0.4089917108642709
$arr[16,0] = 0.2463577109360529
$arr[16,1] = 0.4054108698657046
$arr[17,0] = 0.2459951888141916
$arr[17,1] = 0.4042030948794135
$arr[18,0] = 0.2476347293353456
$arr[18,1] = 0.4096566485440718
$arr[19,0] = 0.2532864529127892
$arr[19,1] = 0.4283037202635285
$arr[20,0] = 0.2570873905531528
$arr[20,1] = 0.4407336843124838
$arr[21,0] = 0.2550489123423318
$arr[21,1] = 0.434077234923933
$arr[22,0] = 0.247544156695227
$arr[22,1] = 0.4093559510764138
$arr[23,0] = 0.2399057176902488
$arr[23,1] = 0.3837228774861359
$ws.Range("L2:M25").Value = $arr

Write-Host "Applied 192 cell updates (rows 2-25, columns B,C,E,F,G,I,L,M)"